# Add a new hourly price-snapshot column to the LDLC smartphone tracking
# sheet. A fresh scrape column is inserted right before the "nom" /
# "url_produit" columns, pushing them one column to the right
# (CA -> CB, CB -> CC). The new column is seeded with:
#   - row 1 (header): the new snapshot's timestamp
#   - rows 2-80: the same price that was already recorded for that row in
#     the previous snapshot column (BZ), since the price did not change
#   - rows 81-206: left blank, matching the blank BZ cells for those rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers: BZ=78, CA=79 (new), CB=80 (was CA), CC=81 (was CB)
$lastSnapshotCol = 78   # BZ
$newCol = 79            # CA after insert

# Insert a new blank column at CA; this shifts the existing "nom" (CA) and
# "url_produit" (CB) columns one position to the right, carrying over their
# values/styles automatically.
$ws.Columns("CA:CA").Insert()

# Header cell for the freshly inserted column.
$ws.Range("CA1").Value2 = "2026-01-31 07:20:34"

# Find the last used row so we cover the whole table.
$lastRow = $ws.UsedRange.Rows.Count

# Carry forward the previous snapshot's price (column BZ) into the new
# column for every data row that had a recorded price.
for ($r = 2; $r -le $lastRow; $r++) {
    $prevValue = $ws.Cells.Item($r, $lastSnapshotCol).Value2
    if ($null -ne $prevValue -and $prevValue -ne "") {
        $ws.Cells.Item($r, $newCol).Value2 = $prevValue
    }
}
